# The "Loại Email" column (originally column L) needs to move to become
# column E, shifting the in-between columns (old E..K) one slot to the
# right. We reproduce this the way a user would in Excel: cut the whole
# column L and insert the cut cells before column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOADON-BACHKHOA")

$ws.Columns("L").Cut() | Out-Null
$ws.Columns("E").Insert() | Out-Null

# Leave the cursor where the author's saved file shows it.
$ws.Range("F12").Select() | Out-Null
